$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the data range keeps text formatting so values are not
# auto-converted to numbers/dates by Excel when assigned below.
$ws.Range("A2:G61").NumberFormat = "@"

$ws.Range("G2").Value = "2025-07-21"
$ws.Range("G3").Value = "2025-07-21"
$ws.Range("B4").Value = "1.39"
$ws.Range("C4").Value = "1.39"
$ws.Range("G4").Value = "2025-07-21"
$ws.Range("G5").Value = "2025-07-21"
$ws.Range("G6").Value = "2025-07-21"
$ws.Range("G7").Value = "2025-07-21"
$ws.Range("G8").Value = "2025-07-21"
$ws.Range("A9").Value = "Bekijk AH Reep melk hazelnoot rozijn"
$ws.Range("B9").Value = "2.59"
$ws.Range("C9").Value = "2.59"
$ws.Range("G9").Value = "2025-07-21"
$ws.Range("A10").Value = "Bekijk AH Praliné bloc wit"
$ws.Range("B10").Value = "2.39"
$ws.Range("C10").Value = "2.39"
$ws.Range("G10").Value = "2025-07-21"
$ws.Range("G11").Value = "2025-07-21"
$ws.Range("A12").Value = "Bekijk AH Speculaas witte chocolade crunchy"
$ws.Range("B12").Value = "1.49"
$ws.Range("C12").Value = "1.49"
$ws.Range("G12").Value = "2025-07-21"
$ws.Range("A13").Value = "Bekijk AH Reep melk hazelnoot rozijn"
$ws.Range("B13").Value = "1.29"
$ws.Range("C13").Value = "1.29"
$ws.Range("D13").Value = "100 g"
$ws.Range("G13").Value = "2025-07-21"
$ws.Range("A14").Value = "Bekijk AH Reep melkchocolade dubbel karamel"
$ws.Range("B14").Value = "2.29"
$ws.Range("C14").Value = "2.29"
$ws.Range("D14").Value = "200 g"
$ws.Range("G14").Value = "2025-07-21"
$ws.Range("G15").Value = "2025-07-21"
$ws.Range("A16").Value = "Bekijk AH Vegan karamel zeezout reep"
$ws.Range("B16").Value = "2.19"
$ws.Range("C16").Value = "2.19"
$ws.Range("D16").Value = "100 g"
$ws.Range("G16").Value = "2025-07-21"
$ws.Range("A17").Value = "Bekijk AH Cookies cream melkchocolade crunchy"
$ws.Range("B17").Value = "1.49"
$ws.Range("C17").Value = "1.49"
$ws.Range("D17").Value = "100 g"
$ws.Range("G17").Value = "2025-07-21"
$ws.Range("A18").Value = "Bekijk AH Reep melkchocolade cookes & cream"
$ws.Range("B18").Value = "2.49"
$ws.Range("C18").Value = "2.49"
$ws.Range("D18").Value = "190 g"
$ws.Range("G18").Value = "2025-07-21"
$ws.Range("A19").Value = "Bekijk AH Reep melkchocolade pinda karamel"
$ws.Range("B19").Value = "2.49"
$ws.Range("C19").Value = "2.49"
$ws.Range("D19").Value = "190 g"
$ws.Range("G19").Value = "2025-07-21"
$ws.Range("G20").Value = "2025-07-21"
$ws.Range("G21").Value = "2025-07-21"
$ws.Range("G22").Value = "2025-07-21"
$ws.Range("G23").Value = "2025-07-21"
$ws.Range("G24").Value = "2025-07-21"
$ws.Range("G25").Value = "2025-07-21"
$ws.Range("G26").Value = "2025-07-21"
$ws.Range("G27").Value = "2025-07-21"
$ws.Range("B28").Value = "1.99"
$ws.Range("C28").Value = "1.99"
$ws.Range("D28").Value = "100 g"
$ws.Range("G28").Value = "2025-07-21"
$ws.Range("B29").Value = "3.59"
$ws.Range("C29").Value = "3.59"
$ws.Range("D29").Value = "200 g"
$ws.Range("G29").Value = "2025-07-21"
$ws.Range("G30").Value = "2025-07-21"
$ws.Range("B31").Value = "1.99"
$ws.Range("C31").Value = "1.99"
$ws.Range("D31").Value = "100 g"
$ws.Range("G31").Value = "2025-07-21"
$ws.Range("B32").Value = "3.59"
$ws.Range("C32").Value = "3.59"
$ws.Range("D32").Value = "200 g"
$ws.Range("G32").Value = "2025-07-21"
$ws.Range("A33").Value = "Bekijk Delicata Reep krachtige pure Belgische chocolade"
$ws.Range("B33").Value = "1.99"
$ws.Range("C33").Value = "1.99"
$ws.Range("D33").Value = "100 g"
$ws.Range("G33").Value = "2025-07-21"
$ws.Range("G34").Value = "2025-07-21"
$ws.Range("A35").Value = "Bekijk Delicata Reep melkchocolade hazelnoot"
$ws.Range("G35").Value = "2025-07-21"
$ws.Range("G36").Value = "2025-07-21"
$ws.Range("G37").Value = "2025-07-21"
$ws.Range("A38").Value = "Bekijk Delicata Reep crispy witte chocolade"
$ws.Range("G38").Value = "2025-07-21"
$ws.Range("A39").Value = "Bekijk Delicata Reep witte chocolade met vanille aroma"
$ws.Range("B39").Value = "3.59"
$ws.Range("C39").Value = "3.59"
$ws.Range("D39").Value = "200 g"
$ws.Range("G39").Value = "2025-07-21"
$ws.Range("G40").Value = "2025-07-21"
$ws.Range("A41").Value = "Bekijk Delicata Reep puur fleur de sel 75% cacao"
$ws.Range("B41").Value = "2.69"
$ws.Range("C41").Value = "2.69"
$ws.Range("D41").Value = "100 g"
$ws.Range("G41").Value = "2025-07-21"
$ws.Range("A42").Value = "Bekijk Delicata Reep pure chocolade hazelnoot"
$ws.Range("B42").Value = "1.99"
$ws.Range("C42").Value = "1.99"
$ws.Range("G42").Value = "2025-07-21"
$ws.Range("A43").Value = "Bekijk Delicata Reep melkchocolade"
$ws.Range("B43").Value = "5.99"
$ws.Range("C43").Value = "5.99"
$ws.Range("D43").Value = "400 g"
$ws.Range("G43").Value = "2025-07-21"
$ws.Range("B44").Value = "3.59"
$ws.Range("C44").Value = "3.59"
$ws.Range("D44").Value = "200 g"
$ws.Range("G44").Value = "2025-07-21"
$ws.Range("G45").Value = "2025-07-21"
$ws.Range("G46").Value = "2025-07-21"
$ws.Range("G47").Value = "2025-07-21"
$ws.Range("G48").Value = "2025-07-21"
$ws.Range("G49").Value = "2025-07-21"
$ws.Range("A50").Value = "Bekijk Delicata Reep pure chocolade sinaasappel amandel"
$ws.Range("G50").Value = "2025-07-21"
$ws.Range("A51").Value = "Bekijk Delicata Reep pure chocolade walnoot vijg"
$ws.Range("G51").Value = "2025-07-21"
$ws.Range("A52").Value = "Bekijk Delicata Reep puur gember limoen"
$ws.Range("B52").Value = "2.89"
$ws.Range("C52").Value = "2.89"
$ws.Range("D52").Value = "90 g"
$ws.Range("G52").Value = "2025-07-21"
$ws.Range("A53").Value = "Bekijk Delicata Reep pure chocolade amandel kokos"
$ws.Range("B53").Value = "3.49"
$ws.Range("C53").Value = "3.49"
$ws.Range("D53").Value = "150 g"
$ws.Range("G53").Value = "2025-07-21"
$ws.Range("G54").Value = "2025-07-21"
$ws.Range("G55").Value = "2025-07-21"
$ws.Range("A56").Value = "Bekijk Delicata Reep pure chocolade balsamico & aarbei"
$ws.Range("B56").Value = "2.89"
$ws.Range("C56").Value = "2.89"
$ws.Range("D56").Value = "90 g"
$ws.Range("G56").Value = "2025-07-21"
$ws.Range("A57").Value = "Bekijk Delicata Reep karamel zeezout cheesecakesmaak"
$ws.Range("G57").Value = "2025-07-21"
$ws.Range("A58").Value = "Bekijk Delicata Reep pure chocolade viooltjes & bosbes"
$ws.Range("G58").Value = "2025-07-21"
$ws.Range("A59").Value = "Bekijk Delicata Reep espresso & karamel pure chocolade"
$ws.Range("G59").Value = "2025-07-21"
$ws.Range("A60").Value = "Bekijk Delicata Reep pure chocolade lavendel & karamel"
$ws.Range("G60").Value = "2025-07-21"
$ws.Range("A61").Value = "Bekijk Delicata Reep pinda pretzel karamel melkchocolade"
$ws.Range("B61").Value = "3.49"
$ws.Range("C61").Value = "3.49"
$ws.Range("D61").Value = "150 g"
$ws.Range("G61").Value = "2025-07-21"
